$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.013731
$ws.Range("H2").Value = 0.041193
$ws.Range("I2").Value = 0.0005742988327511807
$ws.Range("J2").Value = 0.0005742988327511806
$ws.Range("M2").Value = 8.142376000000001
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 0.111802964856
$ws.Range("R2").Value = 1.006226683704
$ws.Range("S2").Value = 0.0001000034559333985
$ws.Range("T2").Value = 0.0001000034559333985

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.013731
$ws.Range("H3").Value = 0.041193
$ws.Range("I3").Value = 0.0005742988327511807
$ws.Range("J3").Value = 0.0005742988327511806
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 0.334217268041
$ws.Range("R3").Value = 3.007955412369
$ws.Range("S3").Value = 0.0002989445036611238
$ws.Range("T3").Value = 0.0002989445036611238

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.013731
$ws.Range("H4").Value = 0.041193
$ws.Range("I4").Value = 0.0005742988327511807
$ws.Range("J4").Value = 0.0005742988327511806
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 0.196040700054
$ws.Range("R4").Value = 1.764366300486
$ws.Range("S4").Value = 0.0001753508731566584
$ws.Range("T4").Value = 0.0001753508731566583

# Row 5
$ws.Range("I5").Value = 0.9843840851703864
$ws.Range("J5").Value = 0.9843840851703864
$ws.Range("M5").Value = 8.142376000000001
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 191.637267921792
$ws.Range("R5").Value = 1724.735411296128
$ws.Range("S5").Value = 0.1714121723202705
$ws.Range("T5").Value = 0.1714121723202705

# Row 6
$ws.Range("I6").Value = 0.9843840851703864
$ws.Range("J6").Value = 0.9843840851703864
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("S6").Value = 0.5124095592244883
$ws.Range("T6").Value = 0.5124095592244883

# Row 7
$ws.Range("I7").Value = 0.9843840851703864
$ws.Range("J7").Value = 0.9843840851703864
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("S7").Value = 0.3005623536256277
$ws.Range("T7").Value = 0.3005623536256277

# Row 8
$ws.Range("I8").Value = 0.01504161599686235
$ws.Range("J8").Value = 0.01504161599686234
$ws.Range("M8").Value = 8.142376000000001
$ws.Range("N8").Value = 24.427128
$ws.Range("O8").Value = 0.1741313933276368
$ws.Range("P8").Value = 0.1741313933276368
$ws.Range("Q8").Value = 2.928261679757334
$ws.Range("R8").Value = 26.354355117816
$ws.Range("S8").Value = 0.00261921755143291
$ws.Range("T8").Value = 0.00261921755143291

# Row 9
$ws.Range("I9").Value = 0.01504161599686235
$ws.Range("J9").Value = 0.01504161599686234
$ws.Range("O9").Value = 0.5205382400466131
$ws.Range("P9").Value = 0.5205382400466131
$ws.Range("S9").Value = 0.007829736318463707
$ws.Range("T9").Value = 0.007829736318463705

# Row 10
$ws.Range("I10").Value = 0.01504161599686235
$ws.Range("J10").Value = 0.01504161599686234
$ws.Range("O10").Value = 0.3053303666257501
$ws.Range("P10").Value = 0.3053303666257501
$ws.Range("S10").Value = 0.004592662126965728
$ws.Range("T10").Value = 0.004592662126965728

